# Generate Report for Handoff
# Adds a new localization-status row (for file
# "2ddbc563-38e3-4348-9953-9a48dd1be2f6.md") to the Overview, zh-cn and
# de-de worksheets, mirroring the existing row for the
# 1d4f5fc2-8ad3-43fc-8060-dc861e1c6a21.md file.

$wb = $excel.ActiveWorkbook

$newFile    = "2ddbc563-38e3-4348-9953-9a48dd1be2f6.md"
$newFilePath= "e2e\" + $newFile
$hyperlinkBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef22f9cfbbfa28c2ace24080c5084e0c1457a132/e2e/" + $newFile

# BGR-encoded equivalent of RGB FF6495ED, used by the workbook's HyperLink
# style, so that re-touched hyperlink cells keep the same visible color.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2:G2").Copy() | Out-Null
$wsOverview.Range("A3:G3").Insert() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $newFilePath
$wsOverview.Range("G3").Value = "2016-09-05 02:44:49"

$linkOv = $wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkBase, "", "", $newFilePath)
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = $hyperlinkColor

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3")) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2:P2").Copy() | Out-Null
$wsZhCn.Range("A3:P3").Insert() | Out-Null

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Range("G3").Value = "2ddbc563-38e3-4348-9953-9a48dd1be2f6.c7213cb916521717e46579afd62ae59aa8c7c853.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 02:44:45"

$linkZh = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkBase, "", "", $newFile)
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3")) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2:P2").Copy() | Out-Null
$wsDeDe.Range("A3:P3").Insert() | Out-Null

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Range("G3").Value = "2ddbc563-38e3-4348-9953-9a48dd1be2f6.c7213cb916521717e46579afd62ae59aa8c7c853.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 02:44:49"

$linkDe = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkBase, "", "", $newFile)
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3")) | Out-Null

$wb.Save()
